$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be placed in the newly created row 2 (top of the data)
$newTopA = -581.0752275067952
$newTopB = 64.53752363065307

# New values to be appended at the bottom (rows 20 and 21)
$newBottom = @(
    @(537.16149053607, 55.61023005004148),
    @(578.7924866790967, 55.64588010560414)
)

# Shift existing data rows (2..18) down to (3..19), working from bottom to top
# so we don't overwrite values before they are read.
for ($r = 18; $r -ge 2; $r--) {
    $destRow = $r + 1
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# Set the new row 2 values (style already matches column A/B since row 2 kept its prior formatting)
$ws.Cells.Item(2, 1).Value = $newTopA
$ws.Cells.Item(2, 2).Value = $newTopB

# Append new rows 19 (new), 20 and 21 with matching style copied from column A (row 2, which has the
# data style) to guarantee rows that previously were empty pick up the correct formatting.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(19, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(20, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(21, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(20, 1).Value = $newBottom[0][0]
$ws.Cells.Item(20, 2).Value = $newBottom[0][1]
$ws.Cells.Item(21, 1).Value = $newBottom[1][0]
$ws.Cells.Item(21, 2).Value = $newBottom[1][1]
